# Generate Report for Handback
# Refresh the handback-status timestamps that the CI run recomputed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 05:02:58"

# --- zh-cn sheet: Correspond Handoff / Handback datetimes for the first row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 05:02:53"
$wsZhCn.Range("K2").Value = "2016-08-26 05:03:16"

# --- de-de sheet: Correspond Handoff / Handback datetimes for the first row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 05:02:58"
$wsDeDe.Range("K2").Value = "2016-08-26 05:03:23"
